# Updates market/profit figures (columns H-N: currentAveragePrice*,
# LevePrice*, LeveProfit*) for a handful of leves across the ALC, ARM,
# BSM, CRP, CUL, GSM and WVR sheets, as produced by the scheduled
# market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 290.75
$ws.Range("I4").Value = 225.63637
$ws.Range("J4").Value = 434
$ws.Range("K4").Value = 225.63637
$ws.Range("L4").Value = 434
$ws.Range("M4").Value = -111.63637
$ws.Range("N4").Value = -662

$ws.Range("H76").Value = 3707123.5
$ws.Range("I76").Value = 4118637.2
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 4118637.2
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -4118322.2
$ws.Range("N76").Value = -4130

$ws.Range("H79").Value = 3707123.5
$ws.Range("I79").Value = 4118637.2
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 4118637.2
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -4117545.2
$ws.Range("N79").Value = -5684

$ws.Range("H113").Value = 3489.9
$ws.Range("I113").Value = 2476.25
$ws.Range("J113").Value = 4165.6665
$ws.Range("K113").Value = 2476.25
$ws.Range("L113").Value = 4165.6665
$ws.Range("M113").Value = 777.75
$ws.Range("N113").Value = -10673.6665

$ws.Range("H116").Value = 8414.944
$ws.Range("I116").Value = 10707.692
$ws.Range("J116").Value = 2453.8
$ws.Range("K116").Value = 10707.692
$ws.Range("L116").Value = 2453.8
$ws.Range("M116").Value = -7265.691999999999
$ws.Range("N116").Value = -9337.799999999999

$ws.Range("H138").Value = 3080.671
$ws.Range("I138").Value = 2044.2162
$ws.Range("J138").Value = 3993.738
$ws.Range("K138").Value = 6132.6486
$ws.Range("L138").Value = 11981.214
$ws.Range("M138").Value = -992.6486000000004
$ws.Range("N138").Value = -22261.214

$ws.Range("H140").Value = 66665
$ws.Range("I140").Value = 30000
$ws.Range("J140").Value = 96663.63
$ws.Range("K140").Value = 30000
$ws.Range("L140").Value = 96663.63
$ws.Range("M140").Value = -24820
$ws.Range("N140").Value = -107023.63

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 596437.75
$ws.Range("I32").Value = 7943.1177
$ws.Range("J32").Value = 7266043.5
$ws.Range("K32").Value = 7943.1177
$ws.Range("L32").Value = 7266043.5
$ws.Range("M32").Value = -7656.1177
$ws.Range("N32").Value = -7266617.5

$ws.Range("H74").Value = 761.7692
$ws.Range("I74").Value = 689.9091
$ws.Range("J74").Value = 1157
$ws.Range("K74").Value = 689.9091
$ws.Range("L74").Value = 1157
$ws.Range("M74").Value = 184.0909
$ws.Range("N74").Value = -2905

$ws.Range("H77").Value = 761.7692
$ws.Range("I77").Value = 689.9091
$ws.Range("J77").Value = 1157
$ws.Range("K77").Value = 3449.5455
$ws.Range("L77").Value = 5785
$ws.Range("M77").Value = 918.4545000000003
$ws.Range("N77").Value = -14521

$ws.Range("H102").Value = 4638.75
$ws.Range("I102").Value = 2222
$ws.Range("J102").Value = 8666.666999999999
$ws.Range("K102").Value = 2222
$ws.Range("L102").Value = 8666.666999999999
$ws.Range("M102").Value = -600
$ws.Range("N102").Value = -11910.667

$ws.Range("H132").Value = 2282.7144
$ws.Range("I132").Value = 1830.5264
$ws.Range("J132").Value = 3844.818
$ws.Range("K132").Value = 5491.5792
$ws.Range("L132").Value = 11534.454
$ws.Range("M132").Value = -2961.5792
$ws.Range("N132").Value = -16594.454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1028.7222
$ws.Range("I134").Value = 931.0769
$ws.Range("J134").Value = 1282.6
$ws.Range("K134").Value = 2793.2307
$ws.Range("L134").Value = 3847.8
$ws.Range("M134").Value = -258.2307000000001
$ws.Range("N134").Value = -8917.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10822.6
$ws.Range("I31").Value = 3866.8965
$ws.Range("J31").Value = 15742.488
$ws.Range("K31").Value = 3866.8965
$ws.Range("L31").Value = 15742.488
$ws.Range("M31").Value = -3571.8965
$ws.Range("N31").Value = -16332.488

$ws.Range("H34").Value = 10822.6
$ws.Range("I34").Value = 3866.8965
$ws.Range("J34").Value = 15742.488
$ws.Range("K34").Value = 3866.8965
$ws.Range("L34").Value = 15742.488
$ws.Range("M34").Value = -3664.8965
$ws.Range("N34").Value = -16146.488

$ws.Range("H62").Value = 3782.318
$ws.Range("I62").Value = 3800.3572
$ws.Range("J62").Value = 3750.75
$ws.Range("K62").Value = 3800.3572
$ws.Range("L62").Value = 3750.75
$ws.Range("M62").Value = -3176.3572
$ws.Range("N62").Value = -4998.75

$ws.Range("H65").Value = 3782.318
$ws.Range("I65").Value = 3800.3572
$ws.Range("J65").Value = 3750.75
$ws.Range("K65").Value = 19001.786
$ws.Range("L65").Value = 18753.75
$ws.Range("M65").Value = -15881.786
$ws.Range("N65").Value = -24993.75

$ws.Range("H132").Value = 2364.818
$ws.Range("I132").Value = 1223.7778
$ws.Range("J132").Value = 7499.5
$ws.Range("K132").Value = 3671.3334
$ws.Range("L132").Value = 22498.5
$ws.Range("M132").Value = -1141.3334
$ws.Range("N132").Value = -27558.5

$ws.Range("H133").Value = 35000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 35000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -40060

$ws.Range("H140").Value = 89900
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 89900
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 89900
$ws.Range("N140").Value = -100260

$ws.Range("H141").Value = 41775.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 41775.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 41775.5
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -52135.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 9260257
$ws.Range("I131").Value = 1280.7142
$ws.Range("J131").Value = 12500899
$ws.Range("K131").Value = 3842.1426
$ws.Range("L131").Value = 37502697
$ws.Range("M131").Value = 1197.8574
$ws.Range("N131").Value = -37512777

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26794904
$ws.Range("I70").Value = 43277384
$ws.Range("J70").Value = 10875
$ws.Range("K70").Value = 43277384
$ws.Range("L70").Value = 10875
$ws.Range("M70").Value = -43277114
$ws.Range("N70").Value = -11415

$ws.Range("H73").Value = 26794904
$ws.Range("I73").Value = 43277384
$ws.Range("J73").Value = 10875
$ws.Range("K73").Value = 43277384
$ws.Range("L73").Value = 10875
$ws.Range("M73").Value = -43276448
$ws.Range("N73").Value = -12747

$ws.Range("H132").Value = 2587.9524
$ws.Range("I132").Value = 1823.125
$ws.Range("J132").Value = 5035.4
$ws.Range("K132").Value = 5469.375
$ws.Range("L132").Value = 15106.2
$ws.Range("M132").Value = -2939.375
$ws.Range("N132").Value = -20166.2

$ws.Range("H140").Value = 89799
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 89799
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 89799
$ws.Range("N140").Value = -100159

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1123.5294
$ws.Range("I136").Value = 706.1111
$ws.Range("J136").Value = 1593.125
$ws.Range("K136").Value = 2118.3333
$ws.Range("L136").Value = 4779.375
$ws.Range("M136").Value = 431.6667000000002
$ws.Range("N136").Value = -9879.375
